$d = $word.ActiveDocument

# Locate the split point inside the run that currently reads:
# "...have been completed for the java portion of the project. Only the android and "
# The edit splits that run right after "portion " (before "of the project."),
# which is where Word drops its _GoBack bookmark to mark the last edit location.
$r = $d.Content
$find = $r.Find
$find.Text = "java portion "
$found = $find.Execute()

if ($found) {
    # Collapse the found range to its end -> the boundary between
    # "...java portion " and "of the project. Only the android and ..."
    $r.Collapse(0)  # wdCollapseEnd

    # Re-adding a bookmark named "_GoBack" moves it here and removes the
    # previous "_GoBack" bookmark (the one that sat before "Keep Track of
    # the Sales of Menu Items"), matching Word's own last-edit tracking.
    $r.Bookmarks.Add("_GoBack")
}

Write-Output "done"
